# CU-1wm69wb update user by id
# Adds three new columns (New-firstName, New-lastName, New-pass) with a second
# test user's data to the "apiTest" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiTest")

# New header cells
$ws.Range("E1").Value = "New-firstName"
$ws.Range("F1").Value = "New-lastName"
$ws.Range("G1").Value = "New-pass"

# New data row for the existing user in row 2
$ws.Range("E2").Value = "Srdjan1"
$ws.Range("F2").Value = "Rados1"
$ws.Range("G2").Value = "Test123@"

# Resize columns to accommodate the new data
$ws.Columns.Item(4).ColumnWidth = 12.584
$ws.Columns.Item(5).ColumnWidth = 17.25
$ws.Columns.Item(6).ColumnWidth = 15.417
$ws.Columns.Item(7).ColumnWidth = 19.25

# Reposition the active cell/selection as left by the editor
$ws.Range("E11").Select() | Out-Null
